{"js": "const body = context.document.body;\n\n// The sentence \"Another while ($true) loop is used to continuously monitor\n// the original path's availability.\" becomes \"Another while ($true) loop\n// continuously monitors the original path's availability.\"\n//\n// This is applied as three small, unambiguous edits so the surrounding\n// run/formatting (Times New Roman, 12pt) is preserved exactly.\n\n// 1) \" loop is used \" -> \" loop \"\nlet r1 = body.search(\" loop is used \", { matchCase: true, matchWholeWord: false });\nr1.load(\"text\");\nawait context.sync();\nr1.items[0].insertText(\" loop \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"to continuously monitor\" -> \"continuously monitors\"\nlet r2 = body.search(\"to continuously monitor\", { matchCase: true, matchWholeWord: false });\nr2.load(\"text\");\nawait context.sync();\nr2.items[0].insertText(\"continuously monitors\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Touch the font of the trailing phrase so it becomes its own run\n//    (matching the source formatting) right after \"continuously monitors\".\nlet r3 = body.search(\" the original path's availability\", { matchCase: true, matchWholeWord: false });\nr3.load(\"text\");\nawait context.sync();\nr3.items[0].font.name = \"Times New Roman\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Another while ($true) loop is used to continuously monitor the original\n# path's availability.\" becomes \"Another while ($true) loop continuously\n# monitors the original path's availability.\"\n#\n# Applied as small, unambiguous Range replacements (rather than a single\n# Find/Replacement.Text swap) so formatting (Times New Roman, 12pt) carries\n# over correctly and so we don't trigger smart-quote autocorrection on the\n# untouched \"path's\" apostrophe.\n\n# 1) \" loop is used \" -> \" loop \"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \" loop is used \"\n$rng1.Find.Execute()\n$rng1.Text = \" loop \"\n\n# 2) \"to continuously monitor\" -> \"continuously monitors\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"to continuously monitor\"\n$rng2.Find.Execute()\n$rng2.Text = \"continuously monitors\"\n\n# 3) Re-assert the font on \"continuously monitors\" so it becomes its own run\n$rng2b = $d.Content\n$rng2b.Find.ClearFormatting()\n$rng2b.Find.Text = \"continuously monitors\"\n$rng2b.Find.Execute()\n$rng2b.Font.Name = \"Times New Roman\"\n\n# 4) Re-assert the font on the trailing phrase so it becomes its own run too\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Text = \" the original path's availability\"\n$rng3.Find.Execute()\n$rng3.Font.Name = \"Times New Roman\"\n"}
